# Horarios 141 - actualizacion 14:11:28 (scrape anterior 13:55:44)
# Aplica los cambios descritos por el diff sobre las 3 hojas del libro:
#   1) LP1912       (sheet1): dimension A1:E199 -> A1:E203, 4 filas nuevas + 3 pares reordenados
#   2) LP1912-215   (sheet2): solo se actualiza el timestamp de cabecera
#   3) 6203-6173    (sheet3): dimension A1:E36 -> A1:E37, 1 fila nueva + 1 par reordenado

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Hoja 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

# Cabecera: hora de actualizacion y total de filas
$ws1.Range("A2").Value = "Última actualización: 14:11:28"
$ws1.Range("A3").Value = "Total filas: 198"

# Dos registros con igual Hora_Scrap/Hora_Llegada intercambian de orden (filas 66/67)
$ws1.Range("A66").Value = "08:28:52"
$ws1.Range("C66").Value = "10_OLMOS"
$ws1.Range("D66").Value = 60
$ws1.Range("A67").Value = "08:11:18"
$ws1.Range("C67").Value = "23_HERNANDEZ"
$ws1.Range("D67").Value = 77

# Mismo caso en filas 141/142 (solo difiere la columna Linea)
$ws1.Range("C141").Value = "23_HERNANDEZ"
$ws1.Range("C142").Value = "27_EL RETIRO"

# Mismo caso en filas 178/179
$ws1.Range("A178").Value = "12:53:26"
$ws1.Range("C178").Value = "11_ETCHEVERRY"
$ws1.Range("D178").Value = 84
$ws1.Range("A179").Value = "12:33:02"
$ws1.Range("C179").Value = "27_EL RETIRO"
$ws1.Range("D179").Value = 104

# Nueva fila scrapeada a las 14:11:28 que llega a las 14:29 -> se inserta
# antes de la fila 182 (14:32) para mantener el orden por Hora_Llegada.
$ws1.Rows(182).Insert()
$ws1.Range("A182").Value = "14:11:28"
$ws1.Range("B182").Value = "14:29"
$ws1.Range("C182").Value = "10_OLMOS"
$ws1.Range("D182").Value = 18
$ws1.Range("E182").Value = "LP1912"

# Nueva fila con llegada 15:36 -> se inserta antes de la fila (ahora) 198 (15:41)
$ws1.Rows(198).Insert()
$ws1.Range("A198").Value = "14:11:28"
$ws1.Range("B198").Value = "15:36"
$ws1.Range("C198").Value = "23_HERNANDEZ"
$ws1.Range("D198").Value = 85
$ws1.Range("E198").Value = "LP1912"

# Dos filas nuevas que se agregan al final de la hoja (llegadas 15:56 y 16:05)
$ws1.Range("A202").Value = "14:11:28"
$ws1.Range("B202").Value = "15:56"
$ws1.Range("C202").Value = "27_EL RETIRO"
$ws1.Range("D202").Value = 105
$ws1.Range("E202").Value = "LP1912"

$ws1.Range("A203").Value = "14:11:28"
$ws1.Range("B203").Value = "16:05"
$ws1.Range("C203").Value = "14_ABASTO"
$ws1.Range("D203").Value = 114
$ws1.Range("E203").Value = "LP1912"

# ---------------------------------------------------------------------------
# Hoja 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 14:11:28"

# ---------------------------------------------------------------------------
# Hoja 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 14:11:28"
$ws3.Range("A3").Value = "Total filas: 32"

# Dos registros con igual Hora_Scrap intercambian de orden (filas 20/21)
$ws3.Range("A20").Value = "08:38:24"
$ws3.Range("C20").Value = "215B_LP-P MOR-1 Y 57"
$ws3.Range("D20").Value = 112
$ws3.Range("A21").Value = "08:52:40"
$ws3.Range("C21").Value = "215A_LA PLATA"
$ws3.Range("D21").Value = 98

# Nueva fila agregada al final de la hoja (llegada 16:02)
$ws3.Range("A37").Value = "14:11:28"
$ws3.Range("B37").Value = "16:02"
$ws3.Range("C37").Value = "215C_LA PLATA"
$ws3.Range("D37").Value = 111
$ws3.Range("E37").Value = "L6203"
